# Fruta / hortaliza, semanal
# Insert two new weekly price observations (dated 2021-11-16, serial 44516)
# right before the existing row 389, shifting the remaining data down by two
# rows (old 389..413 -> new 391..415), matching the canonical diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 389-390; everything below shifts down by 2
# (this also carries the D-column date number format down with the shift).
$ws.Range("A389:A390").EntireRow.Insert()

# New row 389: Provincia de Chacabuco, Primera
$ws.Cells.Item(389, 1).Value = 6
$ws.Cells.Item(389, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(389, 3).Value = "Metropolitana"
$ws.Cells.Item(389, 4).Value = 44516
$ws.Cells.Item(389, 5).Value = 13
$ws.Cells.Item(389, 6).Value = 100112012
$ws.Cells.Item(389, 7).Value = "Espinaca"
$ws.Cells.Item(389, 8).Value = "Sin especificar"
$ws.Cells.Item(389, 9).Value = "Primera"
$ws.Cells.Item(389, 10).Value = 350
$ws.Cells.Item(389, 11).Value = 4500
$ws.Cells.Item(389, 12).Value = 5000
$ws.Cells.Item(389, 13).Value = 4729
$ws.Cells.Item(389, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(389, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(389, 16).Value = 473
$ws.Cells.Item(389, 17).Value = 10
$ws.Cells.Item(389, 18).Value = "Hortaliza"

# New row 390: Región Metropolitana, Primera
$ws.Cells.Item(390, 1).Value = 6
$ws.Cells.Item(390, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(390, 3).Value = "Metropolitana"
$ws.Cells.Item(390, 4).Value = 44516
$ws.Cells.Item(390, 5).Value = 13
$ws.Cells.Item(390, 6).Value = 100112012
$ws.Cells.Item(390, 7).Value = "Espinaca"
$ws.Cells.Item(390, 8).Value = "Sin especificar"
$ws.Cells.Item(390, 9).Value = "Primera"
$ws.Cells.Item(390, 10).Value = 420
$ws.Cells.Item(390, 11).Value = 4500
$ws.Cells.Item(390, 12).Value = 5000
$ws.Cells.Item(390, 13).Value = 4714
$ws.Cells.Item(390, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(390, 15).Value = "Región Metropolitana"
$ws.Cells.Item(390, 16).Value = 471
$ws.Cells.Item(390, 17).Value = 10
$ws.Cells.Item(390, 18).Value = "Hortaliza"

# Keep the D-column date number format consistent with the rest of the column
$ws.Range("D389:D390").NumberFormat = "YYYY-MM-DD HH:MM:SS"
